{"js": "// Remove the \"SO\" section: the bold \"SO\" heading, its screenshot image,\n// the red-text error/traceback block, and the trailing bullet about\n// detecting SO with stage/cycle concat. The blank paragraph immediately\n// before \"SO\" and the blank paragraphs after the final bullet are kept.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet startIndex = -1;\nlet endIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (startIndex === -1 && items[i].text === \"SO\") {\n    startIndex = i;\n  }\n  if (items[i].text.indexOf(\"Make sure that it detects SO with stage and cycle\") === 0) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1 || endIndex < startIndex) {\n  throw new Error(\"Could not locate the 'SO' section to remove.\");\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = endIndex; i >= startIndex; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"SO\" section: the bold \"SO\" heading, its screenshot image,\n# the red-text error/traceback block, and the trailing bullet about\n# detecting SO with stage/cycle concat. The blank paragraph immediately\n# before \"SO\" and the blank paragraphs after the final bullet are kept.\n$d = $word.ActiveDocument\n\n$startIdx = -1\n$endIdx = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($startIdx -eq -1 -and $t -eq \"SO\") {\n        $startIdx = $i\n    }\n    if ($t.StartsWith(\"Make sure that it detects SO with stage and cycle\")) {\n        $endIdx = $i\n        break\n    }\n    $i++\n}\n\nif ($startIdx -eq -1 -or $endIdx -eq -1) {\n    throw \"Could not locate the 'SO' section to remove.\"\n}\n\n$startPara = $d.Paragraphs.Item($startIdx)\n$endPara = $d.Paragraphs.Item($endIdx)\n$range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$range.Delete()\n"}
